# Update existing rows 2-21: columns C-H (ax, ay, az, gx, gy, gz) get new sensor readings
# while timestamp (A) and label (B) are unchanged for these rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{}
$data[2] = @(8.406875610351562, -12.23353862762451, 5.69203519821167, -1.558092090690032, 0.5108912702903179, 1.960815749122113)
$data[3] = @(-6.086453437805176, -26.34650039672852, -14.22333812713623, -1.325494741062512, 3.293620617210238, 0.05309254307193068)
$data[4] = @(5.984857559204102, -12.02742099761963, -10.5483455657959, 0.7431222165672173, 1.402724462805504, -1.679901339475376)
$data[5] = @(-9.330942153930664, -13.30614566802978, -3.071574211120605, 2.074821671235908, -1.928939194355179, -2.52905666828155)
$data[6] = @(-9.185368537902832, -8.196459770202637, -3.727138996124268, 3.729947960492444, -3.12722979934469, -2.61889575291606)
$data[7] = @(-11.09085464477539, -10.81818866729736, 13.48275184631348, -1.309740603548675, -3.723056509656808, 1.097793729559807)
$data[8] = @(-5.257553100585938, -7.122329711914063, -5.412375450134277, -10.97659649895218, -5.284385065430541, 1.629254709748432)
$data[9] = @(5.622566699981689, -31.63553428649902, 26.30605506896973, -2.185527533003532, -9.327580903340284, 0.09413627339798614)
$data[10] = @(-31.68562698364257, 8.553699493408203, -25.76754570007324, 4.113174119041981, -6.550736431936665, 2.042246463808018)
$data[11] = @(-25.48157119750977, -8.927703857421875, -22.08453178405762, 6.01964588419901, 0.3379616459596031, -0.06193271655476229)
$data[12] = @(3.429775714874268, -31.35194206237793, 6.574334621429443, 2.943291461004813, 0.339280297455288, -1.968303656115117)
$data[13] = @(-11.71933746337891, -2.056103706359864, 0.026987075805664, -1.721067187566203, 2.170973975681572, 0.1046053223239598)
$data[14] = @(-4.002825736999512, 1.024898052215576, 5.117716789245605, -3.235397556453075, -3.147966375628718, -2.00323488411394)
$data[15] = @(-2.464081764221191, -30.21155166625977, 11.36835384368896, -0.1528169541683106, -4.183527493071799, 0.4817037536102458)
$data[16] = @(18.32395362854004, 2.226233005523682, -8.703231811523438, 3.802701440829669, -4.169063769497901, 2.565167040500834)
$data[17] = @(-23.21679878234864, -14.6329288482666, -21.63663291931152, 3.380828262532809, -6.121044473740666, 0.1218900376732013)
$data[18] = @(-1.97741436958313, -62.5509033203125, -12.70803451538086, -0.8874874531643746, -2.666994254565963, 1.347052811997586)
$data[19] = @(-14.78749465942383, 9.399082183837891, 0.7060952186584473, -9.827318640588476, 1.861103559581952, 7.985667338648988)
$data[20] = @(-45.69917297363281, -85.84348297119141, 16.09267616271973, -2.809762274177507, 2.697998086225652, 0.7576566992455813)
$data[21] = @(21.77816009521484, 1.36038064956665, -24.16020965576172, 4.037018136492105, -2.376052915471196, -2.018577781695872)
$data[22] = @(-22.23860740661621, -3.703294038772583, -13.579345703125, 5.506457075331898, -2.59009868774606, -0.0266988300582125)
$data[23] = @(-28.10436058044434, -46.36766815185547, 13.63930892944336, 2.450691733545613, 1.12102667965644, -1.359096467929008)
$data[24] = @(0.2610459327697754, 2.697724103927612, -2.434926986694336, -3.604418606434045, 0.8684271347926146, -0.2955643468111461)
$data[25] = @(8.421352386474609, -4.679108619689941, 5.642993927001953, -3.740167403684082, -2.81262268075148, -0.7813068606610404)
$data[26] = @(12.62557125091553, -29.58051681518555, -5.874727249145508, -0.1918180868463897, 4.700260833629038, -2.898207216586866)
$data[27] = @(17.53614234924316, 6.259291648864746, -4.211094856262207, 3.375462932494072, 2.835833433762467, -3.255329332305403)
$data[28] = @(-11.21247386932373, -9.254820823669434, -18.96420669555664, 3.974689759097054, -8.058607795863464, -3.702626524619664)
$data[29] = @(0.2897729873657226, -15.05545616149902, 21.06244659423828, 1.054417246753779, -4.126515001348001, -1.594708632497015)
$data[30] = @(-14.59784603118896, 12.85461044311523, 1.184048175811768, -4.560519280942842, 2.439459521793583, 8.084467855471887)
$data[31] = @(-74.86586761474609, -84.31977081298828, 76.46285247802734, -1.622863769531381, 3.112303316593201, 3.16303658485438)

for ($r = 2; $r -le 21; $r++) {
    $vals = $data[$r]
    $ws.Cells.Item($r, 3).Value = $vals[0]
    $ws.Cells.Item($r, 4).Value = $vals[1]
    $ws.Cells.Item($r, 5).Value = $vals[2]
    $ws.Cells.Item($r, 6).Value = $vals[3]
    $ws.Cells.Item($r, 7).Value = $vals[4]
    $ws.Cells.Item($r, 8).Value = $vals[5]
}

# Append new rows 22-31 (timestamps 2000..2900) with full data including
# timestamp (A) and label (B).
$timestamps = @{}
$timestamps[22] = 2000
$timestamps[23] = 2100
$timestamps[24] = 2200
$timestamps[25] = 2300
$timestamps[26] = 2400
$timestamps[27] = 2500
$timestamps[28] = 2600
$timestamps[29] = 2700
$timestamps[30] = 2800
$timestamps[31] = 2900

for ($r = 22; $r -le 31; $r++) {
    $vals = $data[$r]
    $ws.Cells.Item($r, 1).Value = $timestamps[$r]
    $ws.Cells.Item($r, 2).Value = "walkingToRunning"
    $ws.Cells.Item($r, 3).Value = $vals[0]
    $ws.Cells.Item($r, 4).Value = $vals[1]
    $ws.Cells.Item($r, 5).Value = $vals[2]
    $ws.Cells.Item($r, 6).Value = $vals[3]
    $ws.Cells.Item($r, 7).Value = $vals[4]
    $ws.Cells.Item($r, 8).Value = $vals[5]
}

Write-Output "done"
